# Set the "h" (float16 / half) row's table-cell background to match the
# other floating point rows ("d" float32, "D" float64) -> FEED99
# (commit: "Set h background to match float type")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# FEED99 as a VBA-style RGB long (0xBBGGRR little-endian packing used by
# Shape.Fill.ForeColor.RGB): R=0xFE, G=0xED, B=0x99
$targetRgb = 0xFE + (0xED * 256) + (0x99 * 65536)

# Locate the "BJData" type-marker table that contains the "h" (float16)
# row -- walk every shape on the slide looking for a table whose first
# column holds the single-character type markers, then find the row
# whose marker cell text is "h".
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            $marker = $tbl.Cell($r, 1).Shape.TextFrame.TextRange.Text
            # Use .Equals() (case-sensitive, ordinal) rather than -eq/-ceq
            # since lowercase "h" (float16/half) must not match uppercase
            # "H" (huge-number) in the other table on this slide.
            if ($marker.Equals("h")) {
                for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                    $tbl.Cell($r, $c).Shape.Fill.ForeColor.RGB = $targetRgb
                }
            }
        }
    }
}
